# edit.ps1 -- apply homework_2.docx revision
# Summary of changes (per the provided diff):
#  1. Truncate the Micro/Macro/Weighted explanation paragraphs, dropping their
#     second sentence but keeping the trailing space (xml:space=preserve).
#  2. Delete the "These values measure impurity..." paragraph.
#  3. Delete the "A higher information gain..." paragraph.
#  4. Insert a new, otherwise-empty paragraph (whose paragraph mark carries an
#     eastAsia rFonts hint + lang) right before the "3. Problem 3 Answers:" heading.
#  5. Move the "Using the first two principal components..." block (and the
#     True/False positive rate figures + "Adding more components..." sentence)
#     so it sits right after "TPR: 0.94" (instead of after the old bookmark
#     paragraph), drop the bookmark from the "TPR: 0.94" paragraph, replace the
#     "Yes, by using continuous data..." / "The full dataset provides..."
#     paragraphs with a single "No." paragraph, and move the _GoBack bookmark to
#     its own paragraph that becomes the very last paragraph of the body.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

# -----------------------------------------------------------------------
# Step 1: truncate the Micro / Macro / Weighted sentences.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "Micro: Computes global TP, FP, and FN across all classes before calculating the score. It is useful when considering overall model performance.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Micro: Computes global TP, FP, and FN across all classes before calculating the score. ", 2) | Out-Null

$d.Content.Find.Execute(
    "Macro: Computes the score for each class independently and then averages them. It treats all classes equally, regardless of their distribution.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Macro: Computes the score for each class independently and then averages them. ", 2) | Out-Null

$d.Content.Find.Execute(
    "Weighted: Computes the score for each class and averages them based on class frequency. This accounts for class imbalance.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Weighted: Computes the score for each class and averages them based on class frequency. ", 2) | Out-Null

# -----------------------------------------------------------------------
# Step 2: delete the "These values measure impurity..." paragraph entirely.
# -----------------------------------------------------------------------
$rImpurity = $d.Content
$rImpurity.Find.Execute("These values measure impurity; lower values indicate purer splits.") | Out-Null
$rImpurity.Expand(4) | Out-Null
$rImpurity.Delete()

# -----------------------------------------------------------------------
# Step 3: delete the "A higher information gain..." paragraph entirely.
# -----------------------------------------------------------------------
$rInfoGain = $d.Content
$rInfoGain.Find.Execute("A higher information gain means the feature effectively reduces uncertainty in classification.") | Out-Null
$rInfoGain.Expand(4) | Out-Null
$rInfoGain.Delete()

# -----------------------------------------------------------------------
# Step 4: insert a new, blank paragraph (its mark carries the eastAsia rFonts
# hint + lang) right before "3.Problem 3 Answers:".
# -----------------------------------------------------------------------
$idxProblem3 = Find-ParagraphIndex $d "3.Problem 3 Answers:"
$d.Paragraphs.Item($idxProblem3).Range.InsertParagraphBefore()
$newBlankPara = $d.Paragraphs.Item($idxProblem3)
$newBlankPara.Range.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# -----------------------------------------------------------------------
# Step 5: move the "Using the first two principal components" block (FP/TP/
# FPR/TPR + "Adding more components..."), drop the _GoBack bookmark from the
# "TPR: 0.94" paragraph, collapse the "Yes, by using continuous data..." /
# "The full dataset provides..." paragraphs into a single "No." paragraph,
# and relocate the _GoBack bookmark into its own trailing paragraph.
# -----------------------------------------------------------------------
$idxTprStart = Find-ParagraphIndex $d "TPR: 0.94"
$idxPcaEnd = Find-ParagraphIndex $d "PCA reduces dimensionality but may discard some information, slightly reducing performance compared to using all features."
$idxTrailingBlank = $idxPcaEnd + 1

$pStart = $d.Paragraphs.Item($idxTprStart)
$pEnd = $d.Paragraphs.Item($idxTrailingBlank)
$rBlock = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rBlock.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>TPR: 0.94</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Using the first two principal components: </w:t></w:r></w:p><w:p><w:r><w:t>FP: 2</w:t></w:r></w:p><w:p><w:r><w:t>TP: 47</w:t></w:r></w:p><w:p><w:r><w:t>FPR: 0.05</w:t></w:r></w:p><w:p><w:r><w:t>TPR: 0.96</w:t></w:r></w:p><w:p><w:r><w:t>Adding more components reduces false positives and improves the true positive rate.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>（</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>）</w:t></w:r><w:r><w:t>Does using continuous data benefit the model? How?</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>No.</w:t></w:r></w:p><w:p><w:r><w:t>PCA reduces dimensionality but may discard some information, slightly reducing performance compared to using all features.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
